$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.145.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.274.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.73"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.34"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.53%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.399"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.841.35"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "66.180.34"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.34"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000162"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.242.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "433.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.52"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.06"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.60"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.04%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.417.33"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.503"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.51%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.82"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.93"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.14"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.54"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.78"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.41"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.33"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.754.14"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.773"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.28"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0653"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "321.14"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.14"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0265"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.03%  "
